$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / summary field updates ---
$ws.Range("D5").Value = 'Report Generated On: 08/18/2025 09:47 PM'
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 96

# --- Zero out pricing (H column) for Wed/Thu/Fri/Sat sections (rows 16-73) ---
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("H24").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("H33").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("H36").Value = 0
$ws.Range("H37").Value = 0
$ws.Range("H38").Value = 0
$ws.Range("H39").Value = 0
$ws.Range("H40").Value = 0
$ws.Range("H41").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("H43").Value = 0
$ws.Range("H44").Value = 0
$ws.Range("H45").Value = 0
$ws.Range("H46").Value = 0
$ws.Range("H47").Value = 0
$ws.Range("H48").Value = 0
$ws.Range("H49").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("H52").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("H54").Value = 0
$ws.Range("H55").Value = 0
$ws.Range("H56").Value = 0
$ws.Range("H57").Value = 0
$ws.Range("H58").Value = 0
$ws.Range("H59").Value = 0
$ws.Range("H60").Value = 0
$ws.Range("H61").Value = 0
$ws.Range("H62").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("H68").Value = 0
$ws.Range("H69").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("H72").Value = 0
$ws.Range("H73").Value = 0

# --- Sunday section (rows 78-127): a new record (Point 01 / INS-15-D-S-C) is
#     inserted at the top, pushing every later record down by one row; the
#     TOTAL row moves from 126 to 127, and all pricing in the section becomes 0 ---

# Row 126 currently holds the merged TOTAL label (A126:G126) - copy its cell
# formatting to the new TOTAL row (127) before we unmerge / overwrite it.
$ws.Range("A126:H126").Copy()
$ws.Range("A127:H127").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# The old A126:G126 merge must be removed before row 126 can hold per-column
# data again; the merge then re-appears, shifted, on the new TOTAL row.
$ws.Range("A126:G126").UnMerge()
$ws.Range("A127:G127").Merge()

$ws.Range("A78").Value = 'Point 01'
$ws.Range("B78").Value = 'INS-15-D-S-C'
$ws.Range("D78").Value = 'INS,15kV,Deadend,Polymer,Corr'
$ws.Range("H78").Value = 0
$ws.Range("A79").Value = 'Point 01'
$ws.Range("B79").Value = 'SAA-DI-10-C'
$ws.Range("D79").Value = 'SAA,Dead End I Bolt,1/0,Corr'
$ws.Range("H79").Value = 0
$ws.Range("A80").Value = 'Point 03'
$ws.Range("B80").Value = 'INS-15-P-S-C'
$ws.Range("D80").Value = 'INS,15kV,Pin,Silicon Polymer,Corr'
$ws.Range("H80").Value = 0
$ws.Range("A81").Value = 'Point 03'
$ws.Range("B81").Value = 'PIN-15-PTP-C'
$ws.Range("D81").Value = 'Pin,15kV,Pole top,Corrosive'
$ws.Range("H81").Value = 0
$ws.Range("A82").Value = 'Point 03'
$ws.Range("B82").Value = 'SAA-3-CV-C'
$ws.Range("D82").Value = 'SAA,3 inch,Clevis,Corr'
$ws.Range("H82").Value = 0
$ws.Range("A83").Value = 'Point 05'
$ws.Range("B83").Value = 'INS-15-P-S-C'
$ws.Range("D83").Value = 'INS,15kV,Pin,Silicon Polymer,Corr'
$ws.Range("H83").Value = 0
$ws.Range("A84").Value = 'Point 05'
$ws.Range("B84").Value = 'PIN-15-PTP-C'
$ws.Range("D84").Value = 'Pin,15kV,Pole top,Corrosive'
$ws.Range("H84").Value = 0
$ws.Range("A85").Value = 'Point 05'
$ws.Range("B85").Value = 'POL-40-2'
$ws.Range("D85").Value = 'Pole,40ft,Class 2'
$ws.Range("H85").Value = 0
$ws.Range("A86").Value = 'Point 05'
$ws.Range("B86").Value = 'SAA-3-CV-C'
$ws.Range("D86").Value = 'SAA,3 inch,Clevis,Corr'
$ws.Range("H86").Value = 0
$ws.Range("A87").Value = 'Point 07'
$ws.Range("B87").Value = 'INS-15-P-S-C'
$ws.Range("D87").Value = 'INS,15kV,Pin,Silicon Polymer,Corr'
$ws.Range("H87").Value = 0
$ws.Range("A88").Value = 'Point 07'
$ws.Range("B88").Value = 'PIN-15-PTP-C'
$ws.Range("D88").Value = 'Pin,15kV,Pole top,Corrosive'
$ws.Range("H88").Value = 0
$ws.Range("A89").Value = 'Point 07'
$ws.Range("B89").Value = 'POL-40-2'
$ws.Range("D89").Value = 'Pole,40ft,Class 2'
$ws.Range("H89").Value = 0
$ws.Range("A90").Value = 'Point 07'
$ws.Range("B90").Value = 'SAA-3-CV-C'
$ws.Range("D90").Value = 'SAA,3 inch,Clevis,Corr'
$ws.Range("H90").Value = 0
$ws.Range("A91").Value = 'Point 09'
$ws.Range("B91").Value = 'INS-15-P-S-C'
$ws.Range("D91").Value = 'INS,15kV,Pin,Silicon Polymer,Corr'
$ws.Range("H91").Value = 0
$ws.Range("A92").Value = 'Point 09'
$ws.Range("B92").Value = 'PIN-15-PTP-C'
$ws.Range("D92").Value = 'Pin,15kV,Pole top,Corrosive'
$ws.Range("H92").Value = 0
$ws.Range("A93").Value = 'Point 09'
$ws.Range("B93").Value = 'POL-40-2'
$ws.Range("D93").Value = 'Pole,40ft,Class 2'
$ws.Range("H93").Value = 0
$ws.Range("A94").Value = 'Point 09'
$ws.Range("B94").Value = 'SAA-3-CV-C'
$ws.Range("D94").Value = 'SAA,3 inch,Clevis,Corr'
$ws.Range("H94").Value = 0
$ws.Range("A95").Value = 'Point 11'
$ws.Range("B95").Value = 'INS-15-P-S-C'
$ws.Range("D95").Value = 'INS,15kV,Pin,Silicon Polymer,Corr'
$ws.Range("H95").Value = 0
$ws.Range("A96").Value = 'Point 11'
$ws.Range("B96").Value = 'PIN-15-PTP-C'
$ws.Range("D96").Value = 'Pin,15kV,Pole top,Corrosive'
$ws.Range("H96").Value = 0
$ws.Range("A97").Value = 'Point 11'
$ws.Range("B97").Value = 'POL-40-2'
$ws.Range("D97").Value = 'Pole,40ft,Class 2'
$ws.Range("H97").Value = 0
$ws.Range("A98").Value = 'Point 11'
$ws.Range("B98").Value = 'SAA-3-CV-C'
$ws.Range("D98").Value = 'SAA,3 inch,Clevis,Corr'
$ws.Range("H98").Value = 0
$ws.Range("A99").Value = 'Point 13'
$ws.Range("B99").Value = 'INS-15-P-S-C'
$ws.Range("D99").Value = 'INS,15kV,Pin,Silicon Polymer,Corr'
$ws.Range("H99").Value = 0
$ws.Range("A100").Value = 'Point 13'
$ws.Range("B100").Value = 'PIN-15-PTP-C'
$ws.Range("D100").Value = 'Pin,15kV,Pole top,Corrosive'
$ws.Range("H100").Value = 0
$ws.Range("A101").Value = 'Point 13'
$ws.Range("B101").Value = 'POL-40-2'
$ws.Range("D101").Value = 'Pole,40ft,Class 2'
$ws.Range("H101").Value = 0
$ws.Range("A102").Value = 'Point 13'
$ws.Range("B102").Value = 'SAA-3-CV-C'
$ws.Range("D102").Value = 'SAA,3 inch,Clevis,Corr'
$ws.Range("H102").Value = 0
$ws.Range("A103").Value = 'Point 15'
$ws.Range("B103").Value = 'INS-15-P-S-C'
$ws.Range("D103").Value = 'INS,15kV,Pin,Silicon Polymer,Corr'
$ws.Range("H103").Value = 0
$ws.Range("A104").Value = 'Point 15'
$ws.Range("B104").Value = 'PIN-15-PTP-C'
$ws.Range("D104").Value = 'Pin,15kV,Pole top,Corrosive'
$ws.Range("H104").Value = 0
$ws.Range("A105").Value = 'Point 15'
$ws.Range("B105").Value = 'POL-40-2'
$ws.Range("D105").Value = 'Pole,40ft,Class 2'
$ws.Range("H105").Value = 0
$ws.Range("A106").Value = 'Point 15'
$ws.Range("B106").Value = 'SAA-3-CV-C'
$ws.Range("D106").Value = 'SAA,3 inch,Clevis,Corr'
$ws.Range("H106").Value = 0
$ws.Range("A107").Value = 'Point 17'
$ws.Range("B107").Value = 'INS-15-P-S-C'
$ws.Range("D107").Value = 'INS,15kV,Pin,Silicon Polymer,Corr'
$ws.Range("H107").Value = 0
$ws.Range("A108").Value = 'Point 17'
$ws.Range("B108").Value = 'PIN-15-PTP-C'
$ws.Range("D108").Value = 'Pin,15kV,Pole top,Corrosive'
$ws.Range("H108").Value = 0
$ws.Range("A109").Value = 'Point 17'
$ws.Range("B109").Value = 'POL-40-2'
$ws.Range("D109").Value = 'Pole,40ft,Class 2'
$ws.Range("H109").Value = 0
$ws.Range("A110").Value = 'Point 17'
$ws.Range("B110").Value = 'SAA-3-CV-C'
$ws.Range("D110").Value = 'SAA,3 inch,Clevis,Corr'
$ws.Range("H110").Value = 0
$ws.Range("A111").Value = 'Point 19'
$ws.Range("B111").Value = 'INS-15-P-S-C'
$ws.Range("D111").Value = 'INS,15kV,Pin,Silicon Polymer,Corr'
$ws.Range("H111").Value = 0
$ws.Range("A112").Value = 'Point 19'
$ws.Range("B112").Value = 'PIN-15-PTP-C'
$ws.Range("D112").Value = 'Pin,15kV,Pole top,Corrosive'
$ws.Range("H112").Value = 0
$ws.Range("A113").Value = 'Point 19'
$ws.Range("B113").Value = 'POL-40-2'
$ws.Range("D113").Value = 'Pole,40ft,Class 2'
$ws.Range("H113").Value = 0
$ws.Range("A114").Value = 'Point 19'
$ws.Range("B114").Value = 'SAA-3-CV-C'
$ws.Range("D114").Value = 'SAA,3 inch,Clevis,Corr'
$ws.Range("H114").Value = 0
$ws.Range("A115").Value = 'Point 21'
$ws.Range("B115").Value = 'INS-15-P-S-C'
$ws.Range("D115").Value = 'INS,15kV,Pin,Silicon Polymer,Corr'
$ws.Range("H115").Value = 0
$ws.Range("A116").Value = 'Point 21'
$ws.Range("B116").Value = 'PIN-15-PTP-C'
$ws.Range("D116").Value = 'Pin,15kV,Pole top,Corrosive'
$ws.Range("H116").Value = 0
$ws.Range("A117").Value = 'Point 21'
$ws.Range("B117").Value = 'POL-40-2'
$ws.Range("D117").Value = 'Pole,40ft,Class 2'
$ws.Range("H117").Value = 0
$ws.Range("A118").Value = 'Point 21'
$ws.Range("B118").Value = 'SAA-3-CV-C'
$ws.Range("D118").Value = 'SAA,3 inch,Clevis,Corr'
$ws.Range("H118").Value = 0
$ws.Range("A119").Value = 'Point 35'
$ws.Range("B119").Value = 'INS-15-P-S-C'
$ws.Range("D119").Value = 'INS,15kV,Pin,Silicon Polymer,Corr'
$ws.Range("H119").Value = 0
$ws.Range("A120").Value = 'Point 35'
$ws.Range("B120").Value = 'PIN-15-PTP-C'
$ws.Range("D120").Value = 'Pin,15kV,Pole top,Corrosive'
$ws.Range("H120").Value = 0
$ws.Range("A121").Value = 'Point 35'
$ws.Range("B121").Value = 'POL-40-2'
$ws.Range("D121").Value = 'Pole,40ft,Class 2'
$ws.Range("H121").Value = 0
$ws.Range("A122").Value = 'Point 35'
$ws.Range("B122").Value = 'SAA-3-CV-C'
$ws.Range("D122").Value = 'SAA,3 inch,Clevis,Corr'
$ws.Range("H122").Value = 0
$ws.Range("A123").Value = 'Point 45'
$ws.Range("B123").Value = 'INS-15-P-S-C'
$ws.Range("D123").Value = 'INS,15kV,Pin,Silicon Polymer,Corr'
$ws.Range("H123").Value = 0
$ws.Range("A124").Value = 'Point 45'
$ws.Range("B124").Value = 'PIN-15-PTP-C'
$ws.Range("D124").Value = 'Pin,15kV,Pole top,Corrosive'
$ws.Range("H124").Value = 0
$ws.Range("A125").Value = 'Point 45'
$ws.Range("B125").Value = 'POL-40-2'
$ws.Range("D125").Value = 'Pole,40ft,Class 2'
$ws.Range("H125").Value = 0
$ws.Range("A126").Value = 'Point 45'
$ws.Range("B126").Value = 'SAA-3-CV-C'
$ws.Range("D126").Value = 'SAA,3 inch,Clevis,Corr'
$ws.Range("H126").Value = 0

# New TOTAL row
$ws.Range("A127").Value = "TOTAL"
$ws.Range("H127").Value = 0

